# "Generate Report for Handback" - a handback just completed for a.md/b.md.
# The file that was previously "ready for handoff" (a.md) is now handed back
# and in sync with en-US, and the file that was already in sync (b.md) is now
# the one "ready for handoff" -- i.e. the File-Name <-> role mapping flips
# between row 2 and row 3 on every sheet, the Status text changes, and the
# "Latest Handback DateTime" timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Rebuild the hyperlinks: row2 now displays "b.md" (still pointing at the
# a.md blob url), row3 now displays "a.md" (still pointing at the b.md blob
# url); row4 (.localization-config) is unchanged.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/a.md", "", "", "b.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/b.md", "", "", "a.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/.localization-config", "", "", ".localization-config")

$ws1.Range("B2").Value = $statusHandedBack
$ws1.Range("C2").Value = $statusHandedBack
$ws1.Range("B3").Value = $statusHandedBack
$ws1.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/a.md", "", "", "b.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8081dcf265505f2d1671492912c3a3819a28b2aa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8a97930f2e90da61e4823aa6ad120ba912c261e1/e2e/a.md", "", "", "a.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ee8af2629a5751b925742465ccc7a0dc0d18e20/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/b.md", "", "", "a.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8081dcf265505f2d1671492912c3a3819a28b2aa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8a97930f2e90da61e4823aa6ad120ba912c261e1/e2e/a.md", "", "", "a.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ee8af2629a5751b925742465ccc7a0dc0d18e20/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/.localization-config", "", "", ".localization-config")

$ws2.Range("B2").Value = $statusHandedBack
$ws2.Range("B3").Value = $statusHandedBack
$ws2.Range("G2").Value = "2016-03-03 12:34:57"
$ws2.Range("G3").Value = "2016-03-03 12:34:57"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/a.md", "", "", "b.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bfa14adce946294cf667697cca09b5919d65b3e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/35c40f1d6ed8c45b8c1eda47fdc263996a53cb7f/e2e/a.md", "", "", "a.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/72012e9bc5b76be52381bca11e6f8bf1b2d5bbbf/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/e2e/b.md", "", "", "a.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bfa14adce946294cf667697cca09b5919d65b3e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/35c40f1d6ed8c45b8c1eda47fdc263996a53cb7f/e2e/a.md", "", "", "a.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/72012e9bc5b76be52381bca11e6f8bf1b2d5bbbf/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b498340b13198178b4ffc224543f02959769096c/.localization-config", "", "", ".localization-config")

$ws3.Range("B2").Value = $statusHandedBack
$ws3.Range("B3").Value = $statusHandedBack
$ws3.Range("G2").Value = "2016-03-03 12:35:19"
$ws3.Range("G3").Value = "2016-03-03 12:35:19"
